# Update "National HE" sheet, row 4 (the household expenditure shares),
# columns B through AJ, with the new June data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("National HE")

$values = @(
    0.00800554476769829,
    0.03718093429243172,
    0.01651556812775156,
    0.03014806460128146,
    0.002456945691160374,
    0.0009680266521925276,
    0.001630252534372018,
    0.003338912214397779,
    0.00613332502259482,
    0.0006783218714613196,
    0.001492421298095608,
    0.004607991906929692,
    0.0006073327254271114,
    0.001036342227098777,
    0.003134823020978324,
    0.001047416843595187,
    0.0008243650161362194,
    0.00008984214355762234,
    0.0006088989897662932,
    0.001489625437887362,
    0.005225410486363665,
    0.00439365509946645,
    0.007688805633126627,
    0.004057752348107797,
    0.005297996715271988,
    0.01002773407981162,
    0.003353085672375054,
    0.0006008942301709348,
    0.002139108671080252,
    0.004786818264616974,
    0.01533516290812533,
    0.009190290925008008,
    0.002261397239270025,
    0.005001381294806605,
    0.003155749554641522
)

# Columns B (2) through AJ (36)
$startCol = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(4, $col).Value = $values[$i]
}
